$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table")

$ws.Range("B3").Value = 2
$ws.Range("B4").Value = 1
